$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-12 Thursday", "2025-06-13 Friday"),
    @("305×9=2745", "753×8=6024"),
    @("174×7=1218", "502×8=4016"),
    @("780×5=3900", "396×5=1980"),
    @("383×6=2298", "167×2=334"),
    @("883×4=3532", "719×5=3595"),
    @("948×5=4740", "512×2=1024"),
    @("392×4=1568", "528×7=3696"),
    @("810×7=5670", "599×2=1198"),
    @("531×5=2655", "463×7=3241"),
    @("559×6=3354", "494×5=2470"),
    @("902×4=3608", "866×4=3464"),
    @("728×6=4368", "720×2=1440"),
    @("788×2=1576", "725×9=6525"),
    @("919×4=3676", "712×6=4272"),
    @("626×6=3756", "487×2=974"),
    @("745×5=3725", "461×9=4149"),
    @("792×4=3168", "871×9=7839"),
    @("942×2=1884", "542×6=3252"),
    @("215×4=860", "600×9=5400"),
    @("572×2=1144", "185×4=740"),
    @("558×2=1116", "472×9=4248"),
    @("641×9=5769", "155×7=1085"),
    @("871×2=1742", "881×5=4405"),
    @("529×7=3703", "950×6=5700"),
    @("749×6=4494", "681×6=4086")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
